# About page documentation updates
# Adds the "Policy Setting (2023 USD)" / "FoPITY" schedule-string block
# (rows 4-7) to the BSfGBP sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BSfGBP")

# Row 4: Policy Setting (2023 USD), value in 2023 USD per MWh
# (each cell set individually -> independent, non-shared formulas, as authored)
$ws.Range("A4").Value = "Policy Setting (2023 USD)"
$ws.Range("F4").Formula = "=F2/About!`$A`$12/1000"
$ws.Range("G4").Formula = "=G2/About!`$A`$12/1000"
$ws.Range("H4").Formula = "=H2/About!`$A`$12/1000"
$ws.Range("I4").Formula = "=I2/About!`$A`$12/1000"
$ws.Range("J4").Formula = "=J2/About!`$A`$12/1000"
$ws.Range("K4").Formula = "=K2/About!`$A`$12/1000"
$ws.Range("L4").Formula = "=L2/About!`$A`$12/1000"
$ws.Range("M4").Formula = "=M2/About!`$A`$12/1000"

# Row 5: FoPITY - normalized fraction of peak incentive that year
$ws.Range("A5").Value = "FoPITY"
$ws.Range("F5:M5").Formula = "=F4/MAX(`$D`$4:`$M`$4)"

# Row 6: concatenated schedule fragments
$ws.Range("C6").Value = "(""Schedule 3"",(2021,0),(2024,0),"
$ws.Range("F6:M6").Formula = "=CONCATENATE(""("",F1,"","",F5,""),"")"
$ws.Range("N6").Value = "(2033,0),(2050,0)),"

# Row 7: full concatenated schedule string
$ws.Range("C7").Formula = "=CONCATENATE(C6,F6,G6,H6,I6,J6,K6,L6,M6,N6)"

# Widen column A so the new policy-setting/schedule labels are readable
$ws.Columns.Item(1).ColumnWidth = 23.65

$ws.Range("C8").Select()
